$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 currently holds "Hepatocellular adenoma HNF 1 alpha mutated - Hyperechoic "
# Add the new YouTube link in column D first, then rename the term, so the
# shared-string table is appended in the same order as the source edit.
$ws.Range("D8").Value = "https://youtu.be/91M82AIMyu0"
$ws.Range("B8").Value = "HNF1α-mutated hepatocellular adenoma - Hyperechoic"

# Reflect the resulting selection state (user ends up with B8 selected).
$ws.Range("B8").Select()
